# Update post last call, with all examples and images
#
# 1) The "datetimeFigureOut" date placeholder (08/06/2020 -> 12/06/2020) is
#    repeated on the slide master and on every slide layout. Walk them all
#    and update the cached text of the date placeholder (msoPlaceholderDate
#    = 16) wherever its text still shows the old cached date.
# 2) On slide 1, the shape "Rectangle : coins arrondis 29" has a run whose
#    text reads "specimenRequirements" that must become "specimenRequested".

$p = $ppt.ActivePresentation

$oldDate = "08/06/2020"
$newDate = "12/06/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 1: specimenRequirements -> specimenRequested.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            $runCount = $para.Runs().Count
            for ($ri = 1; $ri -le $runCount; $ri++) {
                $run = $para.Runs($ri)
                # The last run of a paragraph reports a trailing CR in its
                # .Text, so trim before comparing.
                $runText = $run.Text.TrimEnd("`r")
                if ($runText -eq "specimenRequirements") {
                    $run.Text = "specimenRequested"
                }
            }
        }
    }
}
